$wb = $excel.ActiveWorkbook

# Update the RRBS "Donor ID" values for the Donor 3 / H1-H3 rows (57 -> 44)
$wsRRBS = $wb.Worksheets.Item("RRBS")
$wsRRBS.Range("B8").Value = 44
$wsRRBS.Range("B9").Value = 44
$wsRRBS.Range("B10").Value = 44

# Move the RRBS sheet's selection from C1:C10 to G16
$wsRRBS.Range("G16").Select()

# Make ATAC the active/selected tab (it was RNAseq before)
$wsATAC = $wb.Worksheets.Item("ATAC")
$wsATAC.Activate()
